# Auto-generated cell updates derived from the commit diff.
# Each entry: target cell (Row/Col), new display text (Value), and whether
# the cell's NumberFormat must be forced to Text first -- column D stores
# price figures as text, and values such as "351.66" or "0.999" would
# otherwise be auto-coerced to numbers by plain assignment.
$updates = @(
    @{ Row = 2; Col = 'D'; Value = '52.266.99'; ForceText = $false },
    @{ Row = 2; Col = 'E'; Value = '  +0.80%  '; ForceText = $false },
    @{ Row = 3; Col = 'D'; Value = '2.904.79'; ForceText = $false },
    @{ Row = 3; Col = 'E'; Value = '  +3.40%  '; ForceText = $false },
    @{ Row = 4; Col = 'E'; Value = '  +0.04%  '; ForceText = $false },
    @{ Row = 5; Col = 'D'; Value = '351.66'; ForceText = $true },
    @{ Row = 6; Col = 'D'; Value = '111.58'; ForceText = $true },
    @{ Row = 6; Col = 'E'; Value = '  -0.65%  '; ForceText = $false },
    @{ Row = 7; Col = 'E'; Value = '  -0.16%  '; ForceText = $false },
    @{ Row = 8; Col = 'D'; Value = '0.999'; ForceText = $true },
    @{ Row = 8; Col = 'E'; Value = '  -0.08%  '; ForceText = $false },
    @{ Row = 9; Col = 'D'; Value = '0.626'; ForceText = $true },
    @{ Row = 9; Col = 'E'; Value = '  -0.66%  '; ForceText = $false },
    @{ Row = 10; Col = 'D'; Value = '39.77'; ForceText = $true },
    @{ Row = 10; Col = 'E'; Value = '  -1.45%  '; ForceText = $false },
    @{ Row = 11; Col = 'E'; Value = '  +2.63%  '; ForceText = $false },
    @{ Row = 12; Col = 'E'; Value = '  +0.07%  '; ForceText = $false },
    @{ Row = 13; Col = 'D'; Value = '19.83'; ForceText = $true },
    @{ Row = 13; Col = 'E'; Value = '  -0.64%  '; ForceText = $false },
    @{ Row = 14; Col = 'D'; Value = '7.76'; ForceText = $true },
    @{ Row = 14; Col = 'E'; Value = '  -0.31%  '; ForceText = $false },
    @{ Row = 15; Col = 'D'; Value = '3.361.44'; ForceText = $false },
    @{ Row = 15; Col = 'E'; Value = '  +3.26%  '; ForceText = $false },
    @{ Row = 16; Col = 'E'; Value = '  +6.41%  '; ForceText = $false },
    @{ Row = 17; Col = 'D'; Value = '2.921.22'; ForceText = $false },
    @{ Row = 17; Col = 'E'; Value = '  +3.72%  '; ForceText = $false },
    @{ Row = 18; Col = 'D'; Value = '52.273.82'; ForceText = $false },
    @{ Row = 18; Col = 'E'; Value = '  +0.77%  '; ForceText = $false },
    @{ Row = 19; Col = 'D'; Value = '7.61'; ForceText = $true },
    @{ Row = 19; Col = 'E'; Value = '  -0.46%  '; ForceText = $false },
    @{ Row = 20; Col = 'E'; Value = '  +3.76%  '; ForceText = $false },
    @{ Row = 21; Col = 'D'; Value = '14.12'; ForceText = $true },
    @{ Row = 21; Col = 'E'; Value = '  +3.63%  '; ForceText = $false },
    @{ Row = 22; Col = 'E'; Value = '  -0.32%  '; ForceText = $false },
    @{ Row = 23; Col = 'E'; Value = '  +0.52%  '; ForceText = $false },
    @{ Row = 24; Col = 'D'; Value = '269.24'; ForceText = $true },
    @{ Row = 24; Col = 'E'; Value = '  +0.19%  '; ForceText = $false },
    @{ Row = 25; Col = 'D'; Value = '2.75'; ForceText = $true },
    @{ Row = 25; Col = 'E'; Value = '  -1.13%  '; ForceText = $false },
    @{ Row = 26; Col = 'E'; Value = '  +2.04%  '; ForceText = $false },
    @{ Row = 27; Col = 'E'; Value = '  -0.21%  '; ForceText = $false },
    @{ Row = 28; Col = 'E'; Value = '  +2.33%  '; ForceText = $false },
    @{ Row = 29; Col = 'D'; Value = '10.59'; ForceText = $true },
    @{ Row = 29; Col = 'E'; Value = '  +1.76%  '; ForceText = $false },
    @{ Row = 30; Col = 'B'; Value = 'InjectiveProtocol'; ForceText = $false },
    @{ Row = 30; Col = 'C'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; ForceText = $false },
    @{ Row = 30; Col = 'D'; Value = '37.59'; ForceText = $true },
    @{ Row = 30; Col = 'E'; Value = '  -1.78%  '; ForceText = $false },
    @{ Row = 31; Col = 'B'; Value = 'RenderToken'; ForceText = $false },
    @{ Row = 31; Col = 'C'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false },
    @{ Row = 31; Col = 'D'; Value = '6.30'; ForceText = $true },
    @{ Row = 31; Col = 'E'; Value = '  +11.85%  '; ForceText = $false },
    @{ Row = 32; Col = 'D'; Value = '6.57'; ForceText = $true },
    @{ Row = 32; Col = 'E'; Value = '  +6.34%  '; ForceText = $false },
    @{ Row = 33; Col = 'E'; Value = '  +0.04%  '; ForceText = $false },
    @{ Row = 34; Col = 'D'; Value = '0.0979'; ForceText = $true },
    @{ Row = 34; Col = 'E'; Value = '  +11.46%  '; ForceText = $false },
    @{ Row = 35; Col = 'D'; Value = '53.25'; ForceText = $true },
    @{ Row = 35; Col = 'E'; Value = '  +1.61%  '; ForceText = $false },
    @{ Row = 36; Col = 'E'; Value = '  +0.58%  '; ForceText = $false },
    @{ Row = 37; Col = 'E'; Value = '  -0.06%  '; ForceText = $false },
    @{ Row = 38; Col = 'E'; Value = '  +5.10%  '; ForceText = $false },
    @{ Row = 39; Col = 'D'; Value = '18.71'; ForceText = $true },
    @{ Row = 39; Col = 'E'; Value = '  -0.90%  '; ForceText = $false },
    @{ Row = 40; Col = 'B'; Value = 'ARBITRUM'; ForceText = $false },
    @{ Row = 40; Col = 'C'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false },
    @{ Row = 40; Col = 'D'; Value = '2.06'; ForceText = $true },
    @{ Row = 40; Col = 'E'; Value = '  +2.34%  '; ForceText = $false },
    @{ Row = 41; Col = 'B'; Value = 'Stacks'; ForceText = $false },
    @{ Row = 41; Col = 'C'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; ForceText = $false },
    @{ Row = 41; Col = 'D'; Value = '2.84'; ForceText = $true },
    @{ Row = 41; Col = 'E'; Value = '  +13.18%  '; ForceText = $false },
    @{ Row = 42; Col = 'D'; Value = '23.61'; ForceText = $true },
    @{ Row = 42; Col = 'E'; Value = '  +7.09%  '; ForceText = $false },
    @{ Row = 43; Col = 'E'; Value = '  +0.81%  '; ForceText = $false },
    @{ Row = 44; Col = 'D'; Value = '2.64'; ForceText = $true },
    @{ Row = 44; Col = 'E'; Value = '  +9.43%  '; ForceText = $false },
    @{ Row = 45; Col = 'D'; Value = '120.47'; ForceText = $true },
    @{ Row = 45; Col = 'E'; Value = '  -0.46%  '; ForceText = $false },
    @{ Row = 46; Col = 'D'; Value = '2.19'; ForceText = $true },
    @{ Row = 46; Col = 'E'; Value = '  -0.01%  '; ForceText = $false },
    @{ Row = 47; Col = 'E'; Value = '  +3.83%  '; ForceText = $false },
    @{ Row = 48; Col = 'D'; Value = '2.193.37'; ForceText = $false },
    @{ Row = 48; Col = 'E'; Value = '  +4.05%  '; ForceText = $false },
    @{ Row = 49; Col = 'D'; Value = '0.263'; ForceText = $true },
    @{ Row = 49; Col = 'E'; Value = '  +23.49%  '; ForceText = $false },
    @{ Row = 50; Col = 'D'; Value = '0.0335'; ForceText = $true },
    @{ Row = 50; Col = 'E'; Value = '  +11.45%  '; ForceText = $false },
    @{ Row = 51; Col = 'E'; Value = '  +2.39%  '; ForceText = $false }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    $cell = $ws.Range($u.Col + $u.Row)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
